$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A33").Value = 29
$ws.Range("B33").Value = "Varianz"
$ws.Range("C33").Value = "zb 28, aber vorher a scho sicher"
